# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 (the "Type of document / Definition / Why it is
#    important" table) gets a new table style applied
#    ({BD78E881-E53C-4600-82FA-6F51FAE00EE8} -> {CD0B0AFF-CF0F-4484-BD3B-92D86081D9DA}).
#
# 2) The deck's theme (ppt/theme/theme1.xml, the theme used by the slide
#    master / every slide) is switched from the "Integral" design's
#    "Red Violet" colour scheme over to the standard Office theme colour
#    scheme (font scheme / format scheme are identical between the two
#    themes already bundled with this deck, so only the 12 theme colours
#    actually need to change).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{CD0B0AFF-CF0F-4484-BD3B-92D86081D9DA}")
    }
}

# --- 2) Swap the slide master's theme colours to the Office palette -----
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# index : (theme colour slot, new RGB as 0xBBGGRR for the COM RGB property)
$themeColors.Colors(1).RGB  = 0         # dk1      000000
$themeColors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388   # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407     # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456   # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797  # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477   # folHlink 954F72
